$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '94.180.31'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.55%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.105.33'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.20%  '

$ws.Range('E4').Value = '  +0.34%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '614.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.20%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.13'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.49%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.390'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.50%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.21%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.824'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.87%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.116.03'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.29%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.198'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.78%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000244'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.95%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '93.868.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.46%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.55'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.06%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.39'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.31%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.707.33'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.90%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.158.50'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.81%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.02%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.87'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.44%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.19%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '450.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.98%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000201'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.49%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.99'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.60%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.31'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.52%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.66'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.62%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '86.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.28%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.09'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.46%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.294.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.92%  '

$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.03%  '

$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.259'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +14.32%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.181'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.79%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.125'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.05%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.35'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.80%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.01%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.162'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.08%  '

$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.89'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.10%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.19'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.00%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.91'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.04%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.452'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.27%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '479.39'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.74%  '

$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.95'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.86%  '

$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.28'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.81%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.74'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -9.34%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.31'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.11%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '160.20'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.84%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.693'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.82%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.86'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.09%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.44'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.14%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.32'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.41%  '

Write-Host "Applied all changes"